# Update "想去人数" (interested-people count) figures that changed in the
# latest data refresh.
#
#   Sheet "展览"   (worksheet 1): F4 1411 -> 1415, F5 673 -> 674
#   Sheet "全部类型" (worksheet 4): F4 1411 -> 1415, F6 673 -> 674

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1415
$wsExhibit.Range("F5").Value = 674

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1415
$wsAll.Range("F6").Value = 674
